# Recreate the "openai_demo" worksheet: a small table of sample strings in
# column B, a fixed prompt in C1, and per-row AI-powered QUERY() array
# formulas in column D that ask whether each row describes a "cool box".
# Rows 6-11 also generate their own sample text in column B via the same
# custom function, chained over the growing range above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet identity / view -------------------------------------------------
$ws.Name = "openai_demo"

$excel.ActiveWindow.Zoom = 85

# --- column sizing (approximate best-fit widths from the source workbook) --
$ws.Columns.Item(2).ColumnWidth = 47.333333333333336   # B
$ws.Columns.Item(3).ColumnWidth = 17.166666666666668   # C
$ws.Columns.Item(4).ColumnWidth = 254.83333333333334   # D
$ws.Columns.Item(10).ColumnWidth = 20.666666666666668  # J

# --- page setup --------------------------------------------------------
$ws.PageSetup.Orientation = 1   # xlPortrait

# --- static data: sample rows + prompt --------------------------------
$ws.Range("B2").Value = "This is NOT a cool box"
$ws.Range("B1").Value = "This a cool box"
$ws.Range("B3").Value = "This box is just okay"
$ws.Range("B4").Value = "This box is whacky"
$ws.Range("B5").Value = "This box is SUPER cool"
$ws.Range("C1").Value = "Is this a cool box?"

# --- D1:D5 -- per-row "is this a cool box?" query ----------------------
$ws.Range("D1").FormulaArray = '=_xldudf_QUERY(B1, $C$1)'
$ws.Range("D2").FormulaArray = '=_xldudf_QUERY(B2, $C$1)'
$ws.Range("D3").FormulaArray = '=_xldudf_QUERY(B3, $C$1)'
$ws.Range("D4").FormulaArray = '=_xldudf_QUERY(B4, $C$1)'
$ws.Range("D5").FormulaArray = '=_xldudf_QUERY(B5, $C$1)'

# --- rows 6-11: generate a new sample row from the data above, then ----
# --- query it the same way as the seed rows -----------------------------
$genPrompt = "Based on the above data, generate another sample row. Be creative and try not to repeat rows which already exist. Only respond in the form of the row."

for ($r = 6; $r -le 11; $r++) {
    $prevRow = $r - 1
    $bFormula = '=_xldudf_QUERY($B$1:B' + $prevRow + ', "' + $genPrompt + '")'
    $ws.Range("B$r").FormulaArray = $bFormula

    $dFormula = '=_xldudf_QUERY(B' + $r + ', $C$1)'
    $ws.Range("D$r").FormulaArray = $dFormula
}

# --- stray italic-styled, otherwise empty cell left over on row 10 -----
$ws.Range("E10").Font.Italic = $true

# --- selection left where the original author last clicked -------------
$ws.Range("D22").Select() | Out-Null
